$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N11").Value = 6.45
$ws.Range("G20").Value = 1.75
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 2.5
$ws.Range("L20").Value = 5.5
$ws.Range("Q20").Value = 2.35
$ws.Range("R20").Value = 1.57
$ws.Range("Z20").Value = 13
$ws.Range("AD20").Value = 6.5
$ws.Range("AH20").Value = 23
$ws.Range("AI20").Value = 17
$ws.Range("AK20").Value = 41
$ws.Range("AN20").Value = 3.6
$ws.Range("AO20").Value = 9.5
$ws.Range("AQ20").Value = 34
$ws.Range("AR20").Value = 67
$ws.Range("AW20").Value = 6.5
$ws.Range("AX20").Value = 29
$ws.Range("AZ20").Value = 101
$ws.Range("AH35").Value = 16.5
$ws.Range("AT35").Value = 2.3
$ws.Range("AW35").Value = 5
$ws.Range("BB35").Value = 500
$ws.Range("G36").Value = 2.87
$ws.Range("H36").Value = 3.45
$ws.Range("AA36").Value = 23
$ws.Range("AD36").Value = 6.8
$ws.Range("G42").Value = 2.18
$ws.Range("I42").Value = 3.05
$ws.Range("J42").Value = 2.82
$ws.Range("L42").Value = 3.7
$ws.Range("M42").Value = 1.06
$ws.Range("O42").Value = 1.31
$ws.Range("P42").Value = 3.15
$ws.Range("Q42").Value = 1.93
$ws.Range("U42").Value = 1.72
$ws.Range("V42").Value = 2
$ws.Range("W42").Value = 7.8
$ws.Range("Y42").Value = 8.75
$ws.Range("Z42").Value = 21
$ws.Range("AA42").Value = 17.5
$ws.Range("AB42").Value = 27
$ws.Range("AF42").Value = 60
$ws.Range("AI42").Value = 11
$ws.Range("AJ42").Value = 40
$ws.Range("AK42").Value = 27
$ws.Range("AL42").Value = 35
$ws.Range("AN42").Value = 4.1
$ws.Range("AO42").Value = 11.5
$ws.Range("AP42").Value = 20
$ws.Range("AQ42").Value = 45
$ws.Range("AR42").Value = 80
$ws.Range("AS42").Value = 250
$ws.Range("AW42").Value = 5
$ws.Range("AX42").Value = 17.5
$ws.Range("AY42").Value = 25
$ws.Range("AZ42").Value = 90
$ws.Range("BA42").Value = 120
$ws.Range("BB42").Value = 350
